$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Archives")

# Add the new "MD5 Hash" header in column G, matching the bold header
# style already used by the other header cells (A1:F1).
$ws.Range("G1").Value = "MD5 Hash"
$ws.Range("G1").Font.Bold = $true

# Give the new column a sensible width (matches a 32-char MD5 hash).
$ws.Columns.Item(7).ColumnWidth = 35.1666666666667

# Move the active selection (as recorded in the sheet view) to B5.
$ws.Range("B5").Select() | Out-Null

# Rebuild the AutoFilter so it spans the full used range A1:G37
# (dropping the old one first so the ref is replaced rather than just
# getting a new filter criterion appended to the existing range).
$ws.AutoFilterMode = $false
$ws.Range("A1:G37").AutoFilter() | Out-Null

# Keep the workbook-level _FilterDatabase defined name in sync with the
# new AutoFilter range, as Excel normally does automatically.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Archives!_FilterDatabase") {
        $n.RefersTo = "=Archives!`$A`$1:`$G`$37"
    }
}
